$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Split the paragraph that currently reads "One-way Anova" into:
#      "**** " + "note if you do not enter a number of groups no
#      statistics will be run."   (two separate runs)
#    followed by a brand-new paragraph that reads "One-way Anova".
# -----------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "One-way Anova`r") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # Replace the existing text in place with the first run's text.
    $r.Text = "**** "

    # Insert a paragraph break right after "**** " so the note text
    # lands in its own run once the break is removed again.
    $afterStar = $d.Range($target.Range.End - 1, $target.Range.End - 1)
    $afterStar.InsertParagraphAfter()

    $notePara = $target.Next()
    $notePara.Range.Text = "note if you do not enter a number of groups no statistics will be run."

    # Merge the note paragraph back into the "**** " paragraph by
    # deleting the paragraph mark between them -- this keeps the two
    # pieces of text as distinct runs instead of merging them into one.
    $mark = $d.Range($target.Range.End - 1, $target.Range.End)
    $mark.Delete()

    # Now append a new paragraph after the combined one, containing
    # the "One-way Anova" text that used to live here.
    $endOfCombined = $d.Range($target.Range.End - 1, $target.Range.End - 1)
    $endOfCombined.InsertParagraphAfter()
    $anovaPara = $target.Next()
    $anovaPara.Range.Text = "One-way Anova"
}

# -----------------------------------------------------------------
# 2) Rename bookmark __DdeLink__302_1425697230 -> __DdeLink__303_1425697230
#    while preserving the original bookmarkStart ordering (it comes
#    before __DdeLink__12249_2669968041).
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("__DdeLink__302_1425697230")) {
    $bmOld = $d.Bookmarks.Item("__DdeLink__302_1425697230")
    $rOld = $bmOld.Range

    $bmKeepName = "__DdeLink__12249_2669968041"
    $rKeep = $null
    if ($d.Bookmarks.Exists($bmKeepName)) {
        $bmKeep = $d.Bookmarks.Item($bmKeepName)
        $rKeep = $bmKeep.Range
        $bmKeep.Delete()
    }

    $bmOld.Delete()
    $d.Bookmarks.Add("__DdeLink__303_1425697230", $rOld)

    if ($rKeep -ne $null) {
        $d.Bookmarks.Add($bmKeepName, $rKeep)
    }
}
